# Weekly price-data update: a new observation was inserted as row 373
# (Fecha = 2023-07-28, serial 45135), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 373; existing rows 373..470 shift to 374..471
$ws.Rows.Item(373).EntireRow.Insert()

# Populate the newly inserted row 373 with the new weekly observation
$ws.Cells.Item(373, 1).Value  = 5
$ws.Cells.Item(373, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(373, 3).Value  = "Maule"
$ws.Cells.Item(373, 4).Value  = 45135
$ws.Cells.Item(373, 5).Value  = 7
$ws.Cells.Item(373, 6).Value  = 100112008
$ws.Cells.Item(373, 7).Value  = "Coliflor"
$ws.Cells.Item(373, 8).Value  = "Sin especificar"
$ws.Cells.Item(373, 9).Value  = "Primera"
$ws.Cells.Item(373, 10).Value = 5000
$ws.Cells.Item(373, 11).Value = 500
$ws.Cells.Item(373, 12).Value = 500
$ws.Cells.Item(373, 13).Value = 500
$ws.Cells.Item(373, 14).Value = "$/unidad"
$ws.Cells.Item(373, 15).Value = "Región del Maule"
$ws.Cells.Item(373, 16).Value = 500
$ws.Cells.Item(373, 17).Value = 1
$ws.Cells.Item(373, 18).Value = "Hortaliza"
